# Q3 Update - 2025
# Delete the "Benin" and "Pakistan" rows from the 2024 dataset, refresh the
# refugee/asylum-seeker figures for the remaining 2024 rows, renumber the
# "items" column, and roll the global short-url value forward.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Delete the two rows that were dropped from this quarter's dataset ---
# Row 452 = Benin (coo_id 18)
$ws.Rows.Item(452).Delete()
# Pakistan (coo_id 147) was originally row 470; after the Benin shift it is row 469
$ws.Rows.Item(469).Delete()

# --- 2. Refresh refugee / asylum-seeker figures for the 2024 rows (452-476) ---
$updates = @(
    @{ Row = 452; N = "0";      O = "38262"; P = "0";     Q = "0";      T = "0" },
    @{ Row = 453; N = "140";    O = "86";    P = "0";     Q = "0";      T = "0" },
    @{ Row = 454; N = "140";    O = "11902"; P = "0";     Q = "0";      T = "0" },
    @{ Row = 455; N = "17";     O = "131";   P = "0";     Q = "0";      T = "0" },
    @{ Row = 456; N = "5";      O = "6";     P = "0";     Q = "0";      T = "0" },
    @{ Row = 457; N = "49";     O = "18";    P = "5";     Q = "0";      T = "0" },
    @{ Row = 458; N = "55";     O = "38";    P = "0";     Q = "0";      T = "0" },
    @{ Row = 459; N = "30";     O = "14";    P = "0";     Q = "0";      T = "0" },
    @{ Row = 460; N = "0";      O = "5";     P = "0";     Q = "0";      T = "0" },
    @{ Row = 461; N = "0";      O = "5";     P = "0";     Q = "0";      T = "0" },
    @{ Row = 462; N = "0";      O = "18";    P = "0";     Q = "0";      T = "0" },
    @{ Row = 463; N = "0";      O = "9";     P = "0";     Q = "0";      T = "0" },
    @{ Row = 464; N = "5";      O = "5";     P = "0";     Q = "0";      T = "0" },
    @{ Row = 465; N = "7";      O = "5";     P = "0";     Q = "0";      T = "0" },
    @{ Row = 466; N = "122409"; O = "66";    P = "61";    Q = "0";      T = "0" },
    @{ Row = 467; N = "0";      O = "0";     P = "0";     Q = "507438"; T = "46590" },
    @{ Row = 468; N = "245845"; O = "38";    P = "21571"; Q = "0";      T = "0" },
    @{ Row = 469; N = "13";     O = "0";     P = "0";     Q = "0";      T = "0" },
    @{ Row = 470; N = "0";      O = "7";     P = "0";     Q = "0";      T = "0" },
    @{ Row = 471; N = "24";     O = "16";    P = "0";     Q = "0";      T = "0" },
    @{ Row = 472; N = "19";     O = "14";    P = "0";     Q = "0";      T = "0" },
    @{ Row = 473; N = "1063";   O = "1282";  P = "0";     Q = "0";      T = "0" },
    @{ Row = 474; N = "5";      O = "5";     P = "0";     Q = "0";      T = "0" },
    @{ Row = 475; N = "5";      O = "7";     P = "0";     Q = "0";      T = "0" },
    @{ Row = 476; N = "5";      O = "11";    P = "0";     Q = "0";      T = "0" }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Cells.Item($r, 14).Value = $u.N
    $ws.Cells.Item($r, 15).Value = $u.O
    $ws.Cells.Item($r, 16).Value = $u.P
    $ws.Cells.Item($r, 17).Value = $u.Q
    $ws.Cells.Item($r, 20).Value = $u.T
    # Column D ("items") stays a sequential id matching row - 1
    $ws.Cells.Item($r, 4).Value = [string]($r - 1)
}

# --- 3. Roll the global short-url value forward for every data row ---
$lastRow = 476
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 2).Value = "9lW5Ax"
}
